# Apply weekly price/volume update by re-assigning the D, J, K, L, M, P
# column values across rows 2-27 of the "Bruselas (repollito)" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: D,J,K,L,M,P = 44476, 220, 20000, 22000, 20909, 1394
$ws.Range("D2").Value = 44476
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 20909
$ws.Range("P2").Value = 1394

# Row 3: D,J,K,L,M,P = 44398, 130, 20000, 20000, 20000, 1333
$ws.Range("D3").Value = 44398
$ws.Range("J3").Value = 130
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 20000
$ws.Range("P3").Value = 1333

# Row 4: D,J,K,L,M,P = 44449, 220, 22000, 24000, 23091, 1539
$ws.Range("D4").Value = 44449
$ws.Range("J4").Value = 220
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 24000
$ws.Range("M4").Value = 23091
$ws.Range("P4").Value = 1539

# Row 5: D,J,K,L,M,P = 44446, 150, 22000, 24000, 22667, 1511
$ws.Range("D5").Value = 44446
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 22000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 22667
$ws.Range("P5").Value = 1511

# Row 6: D,J,K,L,M,P = 44741, 250, 18000, 20000, 18800, 1253
$ws.Range("D6").Value = 44741
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 18000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 18800
$ws.Range("P6").Value = 1253

# Row 7: D,J,K,L,M,P = 44406, 400, 20000, 22000, 20850, 1390
$ws.Range("D7").Value = 44406
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 20850
$ws.Range("P7").Value = 1390

# Row 8: D,J,K,L,M,P = 44727, 220, 16000, 18000, 16909, 1127
$ws.Range("D8").Value = 44727
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 16909
$ws.Range("P8").Value = 1127

# Row 9: D,J,K,L,M,P = 44742, 400, 18000, 20000, 18850, 1257
$ws.Range("D9").Value = 44742
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 18850
$ws.Range("P9").Value = 1257

# Row 10: D,J,K,L,M,P = 44755, 230, 16000, 18000, 16783, 1119
$ws.Range("D10").Value = 44755
$ws.Range("J10").Value = 230
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 16783
$ws.Range("P10").Value = 1119

# Row 11: D,J,K,L,M,P = 44400, 130, 24000, 24000, 24000, 1600
$ws.Range("D11").Value = 44400
$ws.Range("J11").Value = 130
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 24000
$ws.Range("P11").Value = 1600

# Row 12: D,J,K,L,M,P = 44722, 150, 18000, 20000, 18933, 1262
$ws.Range("D12").Value = 44722
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 18933
$ws.Range("P12").Value = 1262

# Row 13: D,J,K,L,M,P = 44714, 200, 16000, 17000, 16400, 1093
$ws.Range("D13").Value = 44714
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 16400
$ws.Range("P13").Value = 1093

# Row 14: D,J,K,L,M,P = 44699, 150, 18000, 20000, 18667, 1244
$ws.Range("D14").Value = 44699
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 18667
$ws.Range("P14").Value = 1244

# Row 15: D,J,K,L,M,P = 44396, 130, 22000, 22000, 22000, 1467
$ws.Range("D15").Value = 44396
$ws.Range("J15").Value = 130
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 22000
$ws.Range("M15").Value = 22000
$ws.Range("P15").Value = 1467

# Row 16: D,J,K,L,M,P = 44435, 140, 21000, 23000, 21714, 1448
$ws.Range("D16").Value = 44435
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = 21000
$ws.Range("L16").Value = 23000
$ws.Range("M16").Value = 21714
$ws.Range("P16").Value = 1448

# Row 17: D,J,K,L,M,P = 44749, 220, 18000, 20000, 19091, 1273
$ws.Range("D17").Value = 44749
$ws.Range("J17").Value = 220
$ws.Range("K17").Value = 18000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19091
$ws.Range("P17").Value = 1273

# Row 18: D,J,K,L,M,P = 44754, 300, 17000, 19000, 18133, 1209
$ws.Range("D18").Value = 44754
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 19000
$ws.Range("M18").Value = 18133
$ws.Range("P18").Value = 1209

# Row 19: D,J,K,L,M,P = 44365, 580, 20000, 22000, 21103, 1407
$ws.Range("D19").Value = 44365
$ws.Range("J19").Value = 580
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 21103
$ws.Range("P19").Value = 1407

# Row 20: D,J,K,L,M,P = 44736, 180, 17000, 19000, 17889, 1193
$ws.Range("D20").Value = 44736
$ws.Range("J20").Value = 180
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 19000
$ws.Range("M20").Value = 17889
$ws.Range("P20").Value = 1193

# Row 21: D,J,K,L,M,P = 44392, 220, 23000, 23000, 23000, 1533
$ws.Range("D21").Value = 44392
$ws.Range("J21").Value = 220
$ws.Range("K21").Value = 23000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = 23000
$ws.Range("P21").Value = 1533

# Row 22: D,J,K,L,M,P = 44748, 200, 16000, 17000, 16400, 1093
$ws.Range("D22").Value = 44748
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 16000
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = 16400
$ws.Range("P22").Value = 1093

# Row 23: D,J,K,L,M,P = 44483, 220, 18000, 20000, 18909, 1261
$ws.Range("D23").Value = 44483
$ws.Range("J23").Value = 220
$ws.Range("K23").Value = 18000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 18909
$ws.Range("P23").Value = 1261

# Row 24: D,J,K,L,M,P = 44399, 150, 22000, 22000, 22000, 1467
$ws.Range("D24").Value = 44399
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 22000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 22000
$ws.Range("P24").Value = 1467

# Row 25: D,J,K,L,M,P = 44747, 400, 17000, 19000, 17850, 1190
$ws.Range("D25").Value = 44747
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 19000
$ws.Range("M25").Value = 17850
$ws.Range("P25").Value = 1190

# Row 26: D,J,K,L,M,P = 44391, 160, 20000, 20000, 20000, 1333
$ws.Range("D26").Value = 44391
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = 20000
$ws.Range("P26").Value = 1333

# Row 27: D,J,K,L,M,P = 44453, 280, 20000, 22000, 21286, 1419
$ws.Range("D27").Value = 44453
$ws.Range("J27").Value = 280
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 22000
$ws.Range("M27").Value = 21286
$ws.Range("P27").Value = 1419
